$d = $word.ActiveDocument

# Locate the paragraph that ends the "It is not possible to clear selection
# field ..." bullet (last bullet of the "Low" priority section) and insert a
# brand-new bullet right after it, before the blank spacer paragraph that
# precedes "Medium". The new paragraph inherits the same list formatting
# (ListParagraph style, ilvl=1/numId=3 numbering, justify-both) from the
# paragraph it follows.
$anchorText = "It is not possible to clear selection field"
$inserted = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$anchorText*") {
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = "Expense form is blinking in Chrome and Explorer"
        $inserted = $true
        break
    }
}

# Remove the old "Implement server side processing for dataTables" bullet
# entirely (it used to sit right under "Medium").
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Implement server side processing for dataTables*") {
        $p.Range.Delete()
        break
    }
}
